# Fixed minor BOM errors
# Applies the cell-level corrections found in the commit diff to
# TransmitterBoard.xlsx (Sheet1), via the Excel COM object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (Cap / TMK107B7223KA-T): unit price & subtotal corrected ---
$ws.Range("G2").Value = 0.1
$ws.Range("H2").Value = 0.3

# --- Row 12 (Cap / CL21B104KCFNNNE): unit price & subtotal corrected ---
$ws.Range("G12").Value = 0.0497
$ws.Range("H12").Value = 0.0497

# --- Row 14 (Cap / CL21B473KCCWPNC): unit price & subtotal corrected ---
$ws.Range("G14").Value = 0.02996
$ws.Range("H14").Value = 0.02996

# --- Row 18 (Inductor): add manufacturer part number & pricing ---
$ws.Range("C18").Value = "LQM18PN4R7MFRL"
$ws.Range("G18").Value = 0.3387
$ws.Range("H18").Value = 0.3387

# --- Row 23 (PMOS-2): add manufacturer part number, fix footprint & pricing ---
$ws.Range("C23").Value = "AONR21321"
$ws.Range("E23").Value = "TRANS_AONR21321"
$ws.Range("G23").Value = 0.3988
$ws.Range("H23").Value = 0.3988

# --- Row 26 (Res1 / RK73B2ATTD204J): unit price & subtotal corrected ---
$ws.Range("G26").Value = 0.0998
$ws.Range("H26").Value = 0.1996

# --- Row 40 (P9242-RNDGI8): add footprint ---
$ws.Range("E40").Value = "QFN40P600X600X100-49N"
